$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.414.93"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "2.373.30"
$ws.Range("E3").Value = "  +3.17%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "309.44"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").Value = "105.10"
$ws.Range("E6").Value = "  +4.37%  "

$ws.Range("E7").Value = "  -3.50%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").Value = "53.35"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").Value = "2.740.41"
$ws.Range("E15").Value = "  +3.14%  "

$ws.Range("D16").Value = "15.61"
$ws.Range("E16").Value = "  +4.02%  "

$ws.Range("D17").Value = "2.373.31"
$ws.Range("E17").Value = "  +2.98%  "

$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").Value = "43.368.26"
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").Value = "12.03"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("E21").Value = "  +3.80%  "

$ws.Range("D22").Value = "0.0₃0919"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").Value = "68.27"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").Value = "241.55"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "25.83"
$ws.Range("E28").Value = "  +6.09%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.31"
$ws.Range("E29").Value = "  +8.77%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "36.80"
$ws.Range("E30").Value = "  -4.32%  "

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "9.56"
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "161.42"
$ws.Range("E32").Value = "  -3.48%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.27"
$ws.Range("E33").Value = "  -1.02%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "18.41"
$ws.Range("E35").Value = "  +3.59%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.78"
$ws.Range("E36").Value = "  +13.24%  "

$ws.Range("E37").Value = "  +6.17%  "

$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "3.10"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0743"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  +5.87%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.106"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -1.52%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").Value = "  +17.21%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.77"
$ws.Range("E44").Value = "  +3.55%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.007.46"
$ws.Range("E45").Value = "  +1.91%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0290"
$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.14"
$ws.Range("E47").Value = "  +3.70%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "10.62"
$ws.Range("E48").Value = "  +7.85%  "

$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "58.09"
$ws.Range("E49").Value = "  +4.41%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.576.49"
$ws.Range("E51").Value = "  +1.99%  "
